# Update stats for 2025-09 (row 22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6294
$ws.Range("C22").Value = 993
$ws.Range("D22").Value = 5824171
$ws.Range("E22").Value = 925.352875754687
$ws.Range("F22").Value = 8.349113444654854
$ws.Range("G22").Value = 3.870292887029292
$ws.Range("H22").Value = 26.65756922033074
